$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.28770000000003
$ws.Range("E5").Value = 13.10589999999999
$ws.Range("E9").Value = 13.59800000000001
$ws.Range("E11").Value = 13.3671
$ws.Range("A21").Value = -21.11930000000001
$ws.Range("E21").Value = 12.8785
$ws.Range("A23").Value = -21.29890000000002
$ws.Range("A25").Value = -22.39010000000004
